$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: roll forward report week (Volume/Number + date range) ---
$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# --- Row 14 (Murder): numeric updates ---
$ws.Range("D14").Value = 4
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -100
$ws.Range("J14").Value = 31
$ws.Range("K14").Value = -64.516129032258
$ws.Range("L14").Value = -62.068965517241
$ws.Range("M14").Value = -73.809523809523
$ws.Range("N14").Value = -86.25

# F14 becomes a text placeholder "0" (matches style used by C14/C30 when data is N/A)
$ws.Range("F14").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)  # xlPasteFormats - adopt C14 text style

# --- Row 15 (Rape) ---
$ws.Range("C15").Value = 2
$ws.Range("E15").Value = -33.333333333333
$ws.Range("F15").Value = 13
$ws.Range("H15").Value = 18.181818181818
$ws.Range("I15").Value = 86
$ws.Range("J15").Value = 116
$ws.Range("K15").Value = -25.862068965517
$ws.Range("L15").Value = -19.626168224299
$ws.Range("M15").Value = 2.380952380952
$ws.Range("N15").Value = -59.624413145539

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 28
$ws.Range("D16").Value = 34
$ws.Range("E16").Value = -17.647058823529
$ws.Range("F16").Value = 122
$ws.Range("H16").Value = -18.666666666666
$ws.Range("I16").Value = 782
$ws.Range("J16").Value = 876
$ws.Range("K16").Value = -10.730593607305
$ws.Range("L16").Value = 23.149606299212
$ws.Range("M16").Value = -37.44
$ws.Range("N16").Value = -82.660753880266

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 58
$ws.Range("D17").Value = 52
$ws.Range("E17").Value = 11.538461538461
$ws.Range("F17").Value = 224
$ws.Range("G17").Value = 230
$ws.Range("H17").Value = -2.608695652173
$ws.Range("I17").Value = 1662
$ws.Range("J17").Value = 1642
$ws.Range("K17").Value = 1.218026796589
$ws.Range("L17").Value = 16.223776223776
$ws.Range("M17").Value = 63.743842364532
$ws.Range("N17").Value = -24.864376130198

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 36
$ws.Range("E18").Value = 44
$ws.Range("F18").Value = 120
$ws.Range("G18").Value = 104
$ws.Range("H18").Value = 15.384615384615
$ws.Range("I18").Value = 689
$ws.Range("J18").Value = 667
$ws.Range("K18").Value = 3.298350824587
$ws.Range("L18").Value = 27.829313543599
$ws.Range("M18").Value = -43.478260869565
$ws.Range("N18").Value = -87.308896666052

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 80
$ws.Range("D19").Value = 72
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 311
$ws.Range("G19").Value = 287
$ws.Range("H19").Value = 8.362369337979
$ws.Range("I19").Value = 2104
$ws.Range("J19").Value = 2229
$ws.Range("K19").Value = -5.607895917451
$ws.Range("L19").Value = 43.91244870041
$ws.Range("M19").Value = 29.317762753534
$ws.Range("N19").Value = -57.148676171079

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 35
$ws.Range("D20").Value = 31
$ws.Range("E20").Value = 12.903225806451
$ws.Range("F20").Value = 152
$ws.Range("G20").Value = 127
$ws.Range("H20").Value = 19.685039370078
$ws.Range("I20").Value = 1043
$ws.Range("J20").Value = 935
$ws.Range("K20").Value = 11.550802139037
$ws.Range("L20").Value = 67.147435897435
$ws.Range("M20").Value = 3.574975173783
$ws.Range("N20").Value = -90.742877429661

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 239
$ws.Range("D21").Value = 221
$ws.Range("E21").Value = 8.14479638009
$ws.Range("F21").Value = 942
$ws.Range("G21").Value = 915
$ws.Range("H21").Value = 2.950819672131
$ws.Range("I21").Value = 6377
$ws.Range("J21").Value = 6496
$ws.Range("K21").Value = -1.831896551724
$ws.Range("L21").Value = 32.138416908412
$ws.Range("M21").Value = 2.130044843049
$ws.Range("N21").Value = -77.71915726215

# --- Row 22 (Transit) ---
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 12
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 70
$ws.Range("J22").Value = 64
$ws.Range("K22").Value = 9.375
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = -5.405405405405

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -12.5
$ws.Range("F23").Value = 25
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = 13.636363636363
$ws.Range("I23").Value = 148
$ws.Range("J23").Value = 137
$ws.Range("K23").Value = 8.029197080291
$ws.Range("L23").Value = 15.625
$ws.Range("M23").Value = 62.637362637362

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 174
$ws.Range("D24").Value = 203
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 755
$ws.Range("G24").Value = 798
$ws.Range("H24").Value = -5.388471177944
$ws.Range("I24").Value = 5561
$ws.Range("J24").Value = 5852
$ws.Range("K24").Value = -4.972658920027
$ws.Range("L24").Value = 42.73613963039
$ws.Range("M24").Value = 52.272727272727

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 102
$ws.Range("D25").Value = 87
$ws.Range("E25").Value = 17.241379310344
$ws.Range("F25").Value = 369
$ws.Range("G25").Value = 351
$ws.Range("H25").Value = 5.128205128205
$ws.Range("I25").Value = 2696
$ws.Range("J25").Value = 2426
$ws.Range("K25").Value = 11.129431162407
$ws.Range("L25").Value = 35.137844611528
$ws.Range("M25").Value = -3.403797921891

# --- Row 26 (UCR Rape*) ---
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 42.857142857142
$ws.Range("I26").Value = 154
$ws.Range("J26").Value = 184
$ws.Range("K26").Value = -16.304347826087
$ws.Range("L26").Value = -6.666666666666

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("C27").Value = 13
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = 85.714285714285
$ws.Range("F27").Value = 38
$ws.Range("G27").Value = 36
$ws.Range("H27").Value = 5.555555555555
$ws.Range("I27").Value = 260
$ws.Range("J27").Value = 259
$ws.Range("K27").Value = 0.3861003861
$ws.Range("L27").Value = 13.043478260869

# --- Row 28 (Shooting Vic.) ---
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = -87.5
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = -89.473684210526
$ws.Range("I28").Value = 58
$ws.Range("J28").Value = 116
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -55.725190839694
$ws.Range("M28").Value = -52.066115702479
$ws.Range("N28").Value = -80.068728522336

# --- Row 29 (Shooting Inc.) ---
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = -83.333333333333
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 14
$ws.Range("H29").Value = -85.714285714285
$ws.Range("I29").Value = 42
$ws.Range("J29").Value = 88
$ws.Range("K29").Value = -52.272727272727
$ws.Range("L29").Value = -61.467889908256
$ws.Range("M29").Value = -56.701030927835
$ws.Range("N29").Value = -84.386617100371

# --- Row 30 (Hate Crimes) ---
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -66.666666666666
$ws.Range("J30").Value = 22
$ws.Range("K30").Value = 54.545454545454
$ws.Range("L30").Value = 61.904761904761
